$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New example/classification values, entered in the same order the
# --- original author appears to have typed them (I-column examples first,
# --- then the quality-rule note in G2, then the classification labels in H).
$ws.Range("I2").Value = "Proceso Judicial No. 20332-2021-00141"
$ws.Range("I3").Value = "Tenencia y maltrato de Fauna Silvestre (247 COIP)"
$ws.Range("I4").Value = "Con fecha 28 de marzo de 2021 Director de PNG recibe llamada de guardaparque Tobias Robalino quien se encontraba en el aeropuerto de Baltra y avisa sobre el hallazgo de 185 tortugas terrestres en maleta color rojo. Se presenta denuncia el mismo día en Fiscalía por parte de PNG."
$ws.Range("I5").Value = "NIXON POLO DELGADO"
$ws.Range("I6").Value = "SENTENCIA CONDENATORIA EJECUTORIADA"

$ws.Range("G2").Value = "Las variables de identificación deben cumplir con una extensión y combinación de caracteres específicos"

$ws.Range("H2").Value = "Identificación"
$ws.Range("H3").Value = "Evento"
$ws.Range("H4").Value = "Descripción"
$ws.Range("H5").Value = "Evento"
$ws.Range("H6").Value = "Evento"

# --- "No aplica" filled into the Fuente de Datos column for rows 2,4,5,6
$ws.Range("E2").Value = "No aplica"
$ws.Range("E4").Value = "No aplica"
$ws.Range("E5").Value = "No aplica"
$ws.Range("E6").Value = "No aplica"

# --- Highlight fills: yellow on the still-empty Fuente de Datos/Reglas de
# --- Calidad placeholder cells, green down column F (Relaciones).
$ws.Range("E3").Interior.Color = 65535
$ws.Range("G3").Interior.Color = 65535
$ws.Range("F2:F6").Interior.Color = 5296274
$ws.Range("G4:G6").Interior.Color = 65535

# --- Move the active selection to H7, matching where the author left off.
$ws.Activate() | Out-Null
$ws.Range("H7").Select() | Out-Null
